$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "26.800.15"
$ws.Range("E2").Value = "  -1.69%  "
$ws.Range("D3").Value = "1.549.55"
$ws.Range("E3").Value = "  -1.62%  "
$ws.Range("E4").Value = "  -0.09%  "
Set-TextValue "D5" "205.28"
$ws.Range("E5").Value = "  -1.29%  "
Set-TextValue "D6" "0.481"
$ws.Range("E6").Value = "  -1.89%  "
$ws.Range("E7").Value = "  -0.08%  "
Set-TextValue "D8" "21.43"
$ws.Range("E8").Value = "  -3.72%  "
Set-TextValue "D9" "0.245"
$ws.Range("E9").Value = "  -0.99%  "
$ws.Range("E10").Value = "  -1.63%  "
$ws.Range("E11").Value = "  -1.10%  "
$ws.Range("D12").Value = "1.770.36"
$ws.Range("E12").Value = "  -1.63%  "
$ws.Range("D13").Value = "1.549.87"
$ws.Range("E13").Value = "  -1.94%  "
$ws.Range("E14").Value = "  -2.65%  "
Set-TextValue "D15" "0.512"
$ws.Range("E15").Value = "  -1.55%  "
$ws.Range("D16").Value = "26.776.26"
$ws.Range("E16").Value = "  -1.81%  "
Set-TextValue "D17" "61.01"
$ws.Range("E17").Value = "  -2.30%  "
Set-TextValue "D18" "213.61"
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("D19").Value = "0.0₃0687"
$ws.Range("E19").Value = "  +0.14%  "
Set-TextValue "D20" "7.23"
$ws.Range("E20").Value = "  -1.55%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("E22").Value = "  -1.23%  "
Set-TextValue "D23" "9.02"
$ws.Range("E23").Value = "  -3.96%  "
$ws.Range("E24").Value = "  -0.99%  "
Set-TextValue "D25" "152.98"
$ws.Range("E25").Value = "  +0.74%  "
Set-TextValue "D26" "6.50"
$ws.Range("E26").Value = "  -2.85%  "
Set-TextValue "D27" "14.89"
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  -1.94%  "
$ws.Range("E30").Value = "  -0.70%  "
$ws.Range("E31").Value = "  -3.15%  "
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("D33").Value = "1.352.68"
$ws.Range("E33").Value = "  -3.94%  "
Set-TextValue "D34" "2.90"
$ws.Range("E34").Value = "  -0.60%  "
$ws.Range("E35").Value = "  -3.60%  "
$ws.Range("E36").Value = "  -0.92%  "
Set-TextValue "D37" "0.918"
$ws.Range("E37").Value = "  -2.07%  "
$ws.Range("E38").Value = "  -2.09%  "
Set-TextValue "D39" "0.523"
$ws.Range("E39").Value = "  +1.00%  "
$ws.Range("E40").Value = "  -2.19%  "
$ws.Range("E41").Value = "  -0.06%  "
Set-TextValue "D42" "5.59"
$ws.Range("E42").Value = "  +4.48%  "
Set-TextValue "D43" "0.989"
$ws.Range("E43").Value = "  -1.06%  "
$ws.Range("E44").Value = "  +0.29%  "
$ws.Range("E45").Value = "  -2.94%  "
Set-TextValue "D46" "62.89"
$ws.Range("E46").Value = "  -1.42%  "
Set-TextValue "D47" "2.28"
$ws.Range("E47").Value = "  -2.23%  "
$ws.Range("D48").Value = "1.683.64"
$ws.Range("E48").Value = "  -1.69%  "
Set-TextValue "D49" "85.90"
$ws.Range("E49").Value = "  -0.15%  "
Set-TextValue "D50" "0.0506"
$ws.Range("E50").Value = "  +2.44%  "
$ws.Range("E51").Value = "  -1.73%  "
